# Insert a new price-check column ("2026-02-05 22:16:23") right before the
# existing "nom" / "url_produit" columns (GL/GM -> GM/GN), shifting them one
# column to the right. The new column carries forward the last known price
# (copied from the column immediately to its left, i.e. what used to be the
# last timestamp column) for rows that have numeric price history, and stays
# blank for rows that have none.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suivi")

# Column GL = 194th column. Inserting here pushes the former GL ("nom") to GM
# and the former GM ("url_produit") to GN, matching the diff exactly.
$ws.Columns("GL").Insert()

# New header cell for the freshly inserted column.
$ws.Range("GL1").Value() = "2026-02-05 22:16:23"

# Rows 2-80 carry numeric price history in the previous last-timestamp column
# (now shifted to GK, column 193); carry that same price value forward into
# the newly inserted column (column 194).
for ($r = 2; $r -le 80; $r++) {
    $lastPrice = $ws.Cells.Item($r, 193).Value()
    $ws.Cells.Item($r, 194).Value() = $lastPrice
}
